# Updated cryptos list on Mon Oct  7 19:38:25 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns for each coin row
# with newly scraped figures, and swaps the Monero / dogwifhat rows (41/42)
# to reflect their updated ranking order (including their Price/Volume).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values as plain text (e.g. "63.231.96", "1.00"),
# not numbers. Force the column to Text format before writing so Excel
# doesn't "helpfully" reinterpret numeric-looking strings (like "0.999" or
# "1.00") as real numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.312.69"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "2.444.70"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "571.81"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").Value = "146.72"
$ws.Range("E6").Value = "  +0.44%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +1.02%  "

$ws.Range("D9").Value = "2.439.81"
$ws.Range("E9").Value = "  -0.25%  "

$ws.Range("D10").Value = "0.111"
$ws.Range("E10").Value = "  -0.71%  "

$ws.Range("E11").Value = "  +1.20%  "

$ws.Range("D12").Value = "5.24"
$ws.Range("E12").Value = "  -0.99%  "

$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  -0.21%  "

$ws.Range("D14").Value = "27.03"
$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("D15").Value = "0.0000179"
$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").Value = "2.876.29"

$ws.Range("D17").Value = "63.095.07"
$ws.Range("E17").Value = "  +0.84%  "

$ws.Range("D18").Value = "2.439.06"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").Value = "11.33"
$ws.Range("E19").Value = "  +0.50%  "

$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  +5.57%  "

$ws.Range("D21").Value = "327.59"
$ws.Range("E21").Value = "  +0.97%  "

$ws.Range("D22").Value = "4.19"
$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("D23").Value = "2.10"
$ws.Range("E23").Value = "  +13.18%  "

$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").Value = "65.31"
$ws.Range("E25").Value = "  -2.95%  "

$ws.Range("D26").Value = "619.34"
$ws.Range("E26").Value = "  +5.53%  "

$ws.Range("D27").Value = "8.98"
$ws.Range("E27").Value = "  +3.56%  "

$ws.Range("E28").Value = "  +1.46%  "

$ws.Range("D29").Value = "2.561.22"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  +3.36%  "

$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.23%  "

$ws.Range("D32").Value = "8.29"
$ws.Range("E32").Value = "  -2.00%  "

$ws.Range("E33").Value = "  -4.28%  "

$ws.Range("E34").Value = "  +0.96%  "

$ws.Range("D35").Value = "5.22"
$ws.Range("E35").Value = "  +7.27%  "

$ws.Range("D36").Value = "1.53"
$ws.Range("E36").Value = "  -0.63%  "

$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("D38").Value = "0.380"
$ws.Range("E38").Value = "  -0.78%  "

$ws.Range("D39").Value = "5.43"
$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("D40").Value = "18.75"
$ws.Range("E40").Value = "  -0.42%  "

# Row 41 now shows dogwifhat (previously Monero).
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "2.70"
$ws.Range("E41").Value = "  +9.65%  "

# Row 42 now shows Monero (previously dogwifhat).
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "146.19"
$ws.Range("E42").Value = "  -1.76%  "

$ws.Range("E43").Value = "  -1.48%  "

$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("D45").Value = "41.86"
$ws.Range("E45").Value = "  +0.51%  "

$ws.Range("D46").Value = "148.66"
$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("E47").Value = "  +1.93%  "

$ws.Range("D48").Value = "21.19"
$ws.Range("E48").Value = "  +2.79%  "

$ws.Range("D49").Value = "0.0535"
$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("D50").Value = "0.601"
$ws.Range("E50").Value = "  -0.28%  "

$ws.Range("E51").Value = "  +0.56%  "
